$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift contents: A7 becomes the number 1, A8 becomes "xsd:integer", A9 becomes "data"
$ws.Range("A7").Value = 1
$ws.Range("A8").Value = "xsd:integer"
$ws.Range("A9").Value = "data"

# Update the active cell selection to A8
$ws.Range("A8").Select()

# Narrow column A's width
$ws.Columns("A").ColumnWidth = 8.50510204081633
